$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '28.042.56'
$ws.Range('E2').Value2 = '  -2.09%  '

$ws.Range('D3').Value2 = '1.830.56'
$ws.Range('E3').Value2 = '  -1.11%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value2 = '1.000'
$ws.Range('E4').Value2 = '  -0.26%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '325.48'
$ws.Range('E5').Value2 = '  -3.10%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '1.0000'
$ws.Range('E6').Value2 = '  -0.26%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '0.4632'
$ws.Range('E7').Value2 = '  -0.67%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.3874'
$ws.Range('E8').Value2 = '  -1.09%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.07861'
$ws.Range('E9').Value2 = '  -0.23%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '0.9604'
$ws.Range('E10').Value2 = '  -2.36%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '21.88'
$ws.Range('E11').Value2 = '  -1.63%  '

$ws.Range('D12').Value2 = '1.905.71'
$ws.Range('E12').Value2 = '  +2.39%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '5.665'
$ws.Range('E13').Value2 = '  -3.16%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '6.895'
$ws.Range('E14').Value2 = '  -1.87%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '0.06763'
$ws.Range('E15').Value2 = '  -0.03%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '87.10'
$ws.Range('E16').Value2 = '  -0.68%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '1.001'
$ws.Range('E17').Value2 = '  -0.27%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '0.000009938'
$ws.Range('E18').Value2 = '  -1.90%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '16.63'
$ws.Range('E19').Value2 = '  -2.33%  '

$ws.Range('E20').Value2 = '  -0.04%  '

$ws.Range('D21').Value2 = '28.063.78'
$ws.Range('E21').Value2 = '  -1.98%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '5.308'
$ws.Range('E22').Value2 = '  -2.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '10.99'
$ws.Range('E23').Value2 = '  -2.66%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '2.099'
$ws.Range('E24').Value2 = '  -1.16%  '

$ws.Range('D25').Value2 = '2.127.07'
$ws.Range('E25').Value2 = '  +2.19%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '153.71'
$ws.Range('E26').Value2 = '  +0.12%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '19.14'
$ws.Range('E27').Value2 = '  -1.46%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '5.766'
$ws.Range('E28').Value2 = '  -8.29%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '1.975'
$ws.Range('E29').Value2 = '  -2.14%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '117.23'
$ws.Range('E30').Value2 = '  -0.32%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '0.9407'
$ws.Range('E31').Value2 = '  -3.94%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '0.09251'
$ws.Range('E32').Value2 = '  -2.15%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '5.295'
$ws.Range('E33').Value2 = '  -1.56%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '1.316'
$ws.Range('E34').Value2 = '  -2.77%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '3.315'
$ws.Range('E35').Value2 = '  -5.44%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '0.05875'
$ws.Range('E36').Value2 = '  -4.43%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '0.02145'
$ws.Range('E37').Value2 = '  -2.31%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '1.144'
$ws.Range('E38').Value2 = '  -1.44%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '7.783'
$ws.Range('E39').Value2 = '  +2.69%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '0.5595'
$ws.Range('E40').Value2 = '  -1.88%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '9.896'
$ws.Range('E41').Value2 = '  -2.03%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '0.1762'
$ws.Range('E42').Value2 = '  -1.33%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '1.201'
$ws.Range('E43').Value2 = '  -3.83%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '11.65'
$ws.Range('E44').Value2 = '  -1.28%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '0.5276'
$ws.Range('E45').Value2 = '  -1.94%  '

$ws.Range('B46').Value2 = 'RenderToken'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '2.168'
$ws.Range('E46').Value2 = '  -7.94%  '

$ws.Range('B47').Value2 = 'Cronos'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '0.07023'
$ws.Range('E47').Value2 = '  -1.65%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '1.831'
$ws.Range('E48').Value2 = '  -4.18%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '113.13'
$ws.Range('E49').Value2 = '  -1.19%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '0.9993'
$ws.Range('E50').Value2 = '  -0.30%  '

$ws.Range('E51').Value2 = '  +0.00%  '
